$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Handback status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Update every cell that shares this string so the shared-string table collapses
# back down cleanly (matches the sharedStrings.xml diff: the old "Ready for
# handoff" entry disappears and nothing is left referencing it).
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for both locales ---
$wsZhCn.Range("K2").Value = "2016-08-17 18:48:16"
$wsDeDe.Range("K2").Value = "2016-08-17 18:48:24"

# --- Error Detail cleared now that the handback is in sync ---
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (report layout regenerated) ---
# Overview!E:F and the "Status" column (col 3) on both locale sheets widen
# from 17.216 to 29.978 character units; the "Error Detail" column (col 16)
# on both locale sheets narrows from 40 to 13.747 character units.
$wsOverview.Columns("E:F").ColumnWidth = 29.166666666666668

$wsZhCn.Columns("C:C").ColumnWidth = 29.166666666666668
$wsZhCn.Columns("P:P").ColumnWidth = 12.833333333333332

$wsDeDe.Columns("C:C").ColumnWidth = 29.166666666666668
$wsDeDe.Columns("P:P").ColumnWidth = 12.833333333333332
